$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed in the repull/recalculation
$ws.Range("F4").Value = 0
$ws.Range("F19").Value = 3
$ws.Range("F21").Value = 2
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = 1
$ws.Range("F29").Value = 0
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = -1
$ws.Range("F40").Value = -5
